$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "65.765.83"
$ws.Range("E2").Value = "  -1.15%  "

$ws.Range("D3").Value = "3.523.73"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("E4").Value = "  +0.06%  "

Set-TextValue $ws.Range("D5") "599.22"
$ws.Range("E5").Value = "  -1.15%  "

Set-TextValue $ws.Range("D6") "143.99"
$ws.Range("E6").Value = "  -0.90%  "

$ws.Range("D7").Value = "3.523.42"
$ws.Range("E7").Value = "  -0.92%  "

$ws.Range("E8").Value = "  -0.15%  "

Set-TextValue $ws.Range("D9") "0.499"
$ws.Range("E9").Value = "  +0.65%  "

$ws.Range("E10").Value = "  -0.98%  "

Set-TextValue $ws.Range("D11") "7.80"
$ws.Range("E11").Value = "  -2.24%  "

$ws.Range("D13").Value = "4.127.07"
$ws.Range("E13").Value = "  -0.79%  "

Set-TextValue $ws.Range("D14") "0.0000199"
$ws.Range("E14").Value = "  -4.32%  "

Set-TextValue $ws.Range("D15") "28.70"
$ws.Range("E15").Value = "  -4.20%  "

$ws.Range("D16").Value = "3.534.37"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").Value = "65.778.83"
$ws.Range("E18").Value = "  -1.18%  "

Set-TextValue $ws.Range("D19") "10.89"
$ws.Range("E19").Value = "  -5.43%  "

Set-TextValue $ws.Range("D20") "6.21"
$ws.Range("E20").Value = "  -0.07%  "

Set-TextValue $ws.Range("D21") "14.32"
$ws.Range("E21").Value = "  -4.22%  "

Set-TextValue $ws.Range("D22") "414.19"
$ws.Range("E22").Value = "  -3.93%  "

$ws.Range("E23").Value = "  -2.47%  "

Set-TextValue $ws.Range("D24") "77.13"
$ws.Range("E24").Value = "  -2.79%  "

$ws.Range("D25").Value = "3.668.27"
$ws.Range("E25").Value = "  -0.81%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  -2.72%  "

$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D28") "2.44"
$ws.Range("E28").Value = "  -2.53%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D29") "7.78"
$ws.Range("E29").Value = "  -2.93%  "

Set-TextValue $ws.Range("D30") "8.93"
$ws.Range("E30").Value = "  -2.38%  "

$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("D32").Value = "3.523.23"
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("E33").Value = "  -0.72%  "

$ws.Range("E34").Value = "  -3.86%  "

$ws.Range("E35").Value = "  +0.01%  "

Set-TextValue $ws.Range("D36") "7.54"
$ws.Range("E36").Value = "  -3.64%  "

Set-TextValue $ws.Range("D37") "1.30"
$ws.Range("E37").Value = "  -10.85%  "

Set-TextValue $ws.Range("D38") "176.13"
$ws.Range("E38").Value = "  +1.90%  "

Set-TextValue $ws.Range("D39") "5.28"
$ws.Range("E39").Value = "  -5.87%  "

$ws.Range("E40").Value = "  -8.36%  "

Set-TextValue $ws.Range("D41") "0.0821"
$ws.Range("E41").Value = "  -3.14%  "

Set-TextValue $ws.Range("D42") "5.07"
$ws.Range("E42").Value = "  -2.37%  "

$ws.Range("E43").Value = "  -3.39%  "

Set-TextValue $ws.Range("D44") "45.20"
$ws.Range("E44").Value = "  -1.94%  "

$ws.Range("E45").Value = "  -8.14%  "

$ws.Range("E46").Value = "  +0.06%  "

Set-TextValue $ws.Range("D47") "2.38"
$ws.Range("E47").Value = "  -5.87%  "

Set-TextValue $ws.Range("D48") "7.08"
$ws.Range("E48").Value = "  -1.23%  "

Set-TextValue $ws.Range("D49") "22.61"
$ws.Range("E49").Value = "  -3.18%  "

Set-TextValue $ws.Range("D50") "1.09"
$ws.Range("E50").Value = "  -8.19%  "

Set-TextValue $ws.Range("D51") "23.11"
$ws.Range("E51").Value = "  -8.00%  "
